# Apply the "Add database dump and finalize query scripts" edit:
#  - enter the new "duration bucket" data (rows 34-39) under Sheet1
#  - add a fourth chart (stacked column) plotting that data, mirroring
#    the style of the existing "b)" stacked-column chart
#  - move the active selection to the new data / chart area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New source data: Sheet1!A34:B39 --------------------------------------
$ws.Range("A34").Value = "[0-90]"
$ws.Range("B34").Value = 2524

$ws.Range("A35").Value = "[91-180]"
$ws.Range("B35").Value = 180

$ws.Range("A36").Value = "[181-240]"
$ws.Range("B36").Value = 548

$ws.Range("A37").Value = "[241-300]"
$ws.Range("B37").Value = 349

$ws.Range("A38").Value = "[301-360]"
$ws.Range("B38").Value = 188

$ws.Range("A39").Value = "[361-…]"
$ws.Range("B39").Value = 434

# --- New chart: stacked column chart of the duration buckets --------------
$shp = $ws.Shapes.AddChart2(201, 52)
$chart = $shp.Chart
$chart.SetSourceData($ws.Range("A34:B39"))

$chart.HasLegend = $false
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Zadatak 3. c) Broj pesama prema trajanju[s]."

$shp.Name = "Chart 4"

# --- Move the selection to reflect the area of the new work ---------------
$null = $ws.Range("P38").Select()
